$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 103.14286
$ws.Range("J5").Value = 102
$ws.Range("L5").Value = 102
$ws.Range("N5").Value = -332
$ws.Range("H28").Value = 549.7692
$ws.Range("I28").Value = 587.25
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 587.25
$ws.Range("L28").Value = 100
$ws.Range("M28").Value = -102.25
$ws.Range("N28").Value = -1070
$ws.Range("H53").Value = 2096.4
$ws.Range("I53").Value = 150
$ws.Range("K53").Value = 150
$ws.Range("M53").Value = 487
$ws.Range("H80").Value = 600.3043
$ws.Range("I80").Value = 416.42856
$ws.Range("J80").Value = 886.3333
$ws.Range("K80").Value = 1249.28568
$ws.Range("L80").Value = 2658.9999
$ws.Range("M80").Value = -251.28568
$ws.Range("N80").Value = -4654.9999
$ws.Range("H83").Value = 600.3043
$ws.Range("I83").Value = 416.42856
$ws.Range("J83").Value = 886.3333
$ws.Range("K83").Value = 3747.85704
$ws.Range("L83").Value = 7976.9997
$ws.Range("M83").Value = 1244.14296
$ws.Range("N83").Value = -17960.9997
$ws.Range("H98").Value = 1551.3529
$ws.Range("I98").Value = 1669.3334
$ws.Range("J98").Value = 1268.2
$ws.Range("K98").Value = 1669.3334
$ws.Range("L98").Value = 1268.2
$ws.Range("M98").Value = -171.3334
$ws.Range("N98").Value = -4264.2
$ws.Range("H106").Value = 31997.2
$ws.Range("I106").Value = 32150.615
$ws.Range("K106").Value = 32150.615
$ws.Range("M106").Value = -31519.615
$ws.Range("H111").Value = 916.8333
$ws.Range("I111").Value = 928
$ws.Range("K111").Value = 2784
$ws.Range("M111").Value = 283
$ws.Range("H122").Value = 1551.3529
$ws.Range("I122").Value = 1669.3334
$ws.Range("J122").Value = 1268.2
$ws.Range("K122").Value = 5008.0002
$ws.Range("L122").Value = 3804.6
$ws.Range("M122").Value = -2558.0002
$ws.Range("N122").Value = -8704.6
$ws.Range("H138").Value = 3803.25
$ws.Range("I138").Value = 3123.8276
$ws.Range("J138").Value = 5116.8
$ws.Range("K138").Value = 9371.4828
$ws.Range("L138").Value = 15350.4
$ws.Range("M138").Value = -4231.4828
$ws.Range("N138").Value = -25630.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 28330
$ws.Range("J24").Value = 28330
$ws.Range("L24").Value = 28330
$ws.Range("N24").Value = -29078
$ws.Range("H32").Value = 6820.6743
$ws.Range("I32").Value = 4689.054
$ws.Range("K32").Value = 4689.054
$ws.Range("M32").Value = -4402.054
$ws.Range("H61").Value = 112
$ws.Range("I61").Value = 112
$ws.Range("K61").Value = 112
$ws.Range("M61").Value = 100
$ws.Range("H74").Value = 947.73334
$ws.Range("I74").Value = 947.73334
$ws.Range("K74").Value = 947.73334
$ws.Range("M74").Value = -73.73334
$ws.Range("H77").Value = 947.73334
$ws.Range("I77").Value = 947.73334
$ws.Range("K77").Value = 4738.6667
$ws.Range("M77").Value = -370.6666999999998
$ws.Range("H86").Value = 63999
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 63999
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 63999
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -66371
$ws.Range("H89").Value = 63999
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 63999
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 191997
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -203853
$ws.Range("H96").Value = 9055.333000000001
$ws.Range("J96").Value = 9055.333000000001
$ws.Range("L96").Value = 9055.333000000001
$ws.Range("N96").Value = -14547.333
$ws.Range("H100").Value = 28330
$ws.Range("J100").Value = 28330
$ws.Range("L100").Value = 28330
$ws.Range("N100").Value = -30494
$ws.Range("H132").Value = 1968.1666
$ws.Range("I132").Value = 1968.1666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5904.4998
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3374.4998
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 112
$ws.Range("I136").Value = 112
$ws.Range("K136").Value = 336
$ws.Range("M136").Value = 2214

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -21232
$ws.Range("H94").Value = 1743.5714
$ws.Range("I94").Value = 1534.1666
$ws.Range("K94").Value = 1534.1666
$ws.Range("M94").Value = -1083.1666
$ws.Range("H99").Value = 2518.6956
$ws.Range("I99").Value = 2391.625
$ws.Range("J99").Value = 2586.4666
$ws.Range("K99").Value = 2391.625
$ws.Range("L99").Value = 2586.4666
$ws.Range("M99").Value = -893.625
$ws.Range("N99").Value = -5582.4666
$ws.Range("H107").Value = 3308.5557
$ws.Range("I107").Value = 3151.5
$ws.Range("K107").Value = 3151.5
$ws.Range("M107").Value = -1231.5
$ws.Range("H134").Value = 2366.9048
$ws.Range("I134").Value = 1713.6666
$ws.Range("K134").Value = 5140.9998
$ws.Range("M134").Value = -2605.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4092.5
$ws.Range("I31").Value = 3471.111
$ws.Range("J31").Value = 4891.4287
$ws.Range("K31").Value = 3471.111
$ws.Range("L31").Value = 4891.4287
$ws.Range("M31").Value = -3176.111
$ws.Range("N31").Value = -5481.4287
$ws.Range("H34").Value = 4092.5
$ws.Range("I34").Value = 3471.111
$ws.Range("J34").Value = 4891.4287
$ws.Range("K34").Value = 3471.111
$ws.Range("L34").Value = 4891.4287
$ws.Range("M34").Value = -3269.111
$ws.Range("N34").Value = -5295.4287
$ws.Range("H58").Value = 2071.9023
$ws.Range("I58").Value = 1103.8
$ws.Range("K58").Value = 1103.8
$ws.Range("M58").Value = -900.8
$ws.Range("H92").Value = 50666.668
$ws.Range("J92").Value = 50666.668
$ws.Range("L92").Value = 50666.668
$ws.Range("N92").Value = -55658.668
$ws.Range("H132").Value = 1466.8334
$ws.Range("I132").Value = 1466.8334
$ws.Range("K132").Value = 4400.5002
$ws.Range("M132").Value = -1870.5002
$ws.Range("H136").Value = 2071.9023
$ws.Range("I136").Value = 1103.8
$ws.Range("K136").Value = 3311.4
$ws.Range("M136").Value = -761.3999999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 12715.875
$ws.Range("I38").Value = 25287.5
$ws.Range("J38").Value = 144.25
$ws.Range("K38").Value = 75862.5
$ws.Range("L38").Value = 432.75
$ws.Range("M38").Value = -75515.5
$ws.Range("N38").Value = -1126.75
$ws.Range("H40").Value = 158.33333
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 158.33333
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 633.33332
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -771.33332
$ws.Range("H132").Value = 4291.5
$ws.Range("J132").Value = 3103.6
$ws.Range("L132").Value = 27932.4
$ws.Range("N132").Value = -32992.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6680.24
$ws.Range("J70").Value = 8399.799999999999
$ws.Range("L70").Value = 8399.799999999999
$ws.Range("N70").Value = -8939.799999999999
$ws.Range("H73").Value = 6680.24
$ws.Range("J73").Value = 8399.799999999999
$ws.Range("L73").Value = 8399.799999999999
$ws.Range("N73").Value = -10271.8
$ws.Range("H114").Value = 43000
$ws.Range("J114").Value = 43000
$ws.Range("L114").Value = 43000
$ws.Range("N114").Value = -51678
$ws.Range("H126").Value = 5342.6665
$ws.Range("J126").Value = 5342.6665
$ws.Range("L126").Value = 16027.9995
$ws.Range("N126").Value = -20967.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1791.6666
$ws.Range("I16").Value = 1587.5
$ws.Range("K16").Value = 1587.5
$ws.Range("M16").Value = -1417.5
$ws.Range("H93").Value = 1623.9286
$ws.Range("I93").Value = 1303.75
$ws.Range("J93").Value = 2050.8333
$ws.Range("K93").Value = 1303.75
$ws.Range("L93").Value = 2050.8333
$ws.Range("M93").Value = -55.75
$ws.Range("N93").Value = -4546.8333
$ws.Range("H106").Value = 17999.715
$ws.Range("J106").Value = 17999.715
$ws.Range("L106").Value = 17999.715
$ws.Range("N106").Value = -20523.715
$ws.Range("H132").Value = 4429.6875
$ws.Range("I132").Value = 3075.25
$ws.Range("K132").Value = 9225.75
$ws.Range("M132").Value = -6695.75
$ws.Range("H136").Value = 4798
$ws.Range("I136").Value = 4286.75
$ws.Range("K136").Value = 12860.25
$ws.Range("M136").Value = -10310.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 20000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 20000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 20000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -20616
$ws.Range("H97").Value = 32786
$ws.Range("J97").Value = 32786
$ws.Range("L97").Value = 32786
$ws.Range("N97").Value = -34768
$ws.Range("H132").Value = 2965.0625
$ws.Range("I132").Value = 2298.92
$ws.Range("K132").Value = 6896.76
$ws.Range("M132").Value = -4366.76
$ws.Range("H136").Value = 1808.8
$ws.Range("I136").Value = 1193.8182
$ws.Range("K136").Value = 3581.4546
$ws.Range("M136").Value = -1031.4546
